$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.487.06"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.018.06"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'263.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.05%  "

$ws.Range("E6").Value = "  -1.56%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'56.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.51%  "

$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("E10").Value = "  -3.27%  "

$ws.Range("E11").Value = "  -2.00%  "

$ws.Range("D12").Value = "'14.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.14%  "

$ws.Range("D13").Value = "2.315.01"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("E14").Value = "  -4.40%  "

$ws.Range("D15").Value = "'20.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.67%  "

$ws.Range("E16").Value = "  -3.70%  "

$ws.Range("D17").Value = "2.020.90"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").Value = "37.397.62"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").Value = "'69.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.70%  "

$ws.Range("D20").Value = "0.0₃0843"
$ws.Range("E20").Value = "  -2.29%  "

$ws.Range("D21").Value = "'5.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").Value = "'228.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.74%  "

$ws.Range("D23").Value = "'2.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.83%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("D26").Value = "'165.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "

$ws.Range("D27").Value = "'8.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.46%  "

$ws.Range("D28").Value = "'19.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.43%  "

$ws.Range("D29").Value = "'0.130"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.94%  "

$ws.Range("D30").Value = "'1.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").Value = "'0.0652"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("E33").Value = "  -3.29%  "

$ws.Range("D34").Value = "'4.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.58%  "

$ws.Range("D35").Value = "'2.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").Value = "'1.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("E39").Value = "  -4.67%  "

$ws.Range("D40").Value = "'3.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.07%  "

$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("D42").Value = "'0.0937"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.10%  "

$ws.Range("D43").Value = "'0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "

$ws.Range("D44").Value = "1.392.93"
$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("D45").Value = "'90.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("D46").Value = "'15.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.55%  "

$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.21%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'7.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.26%  "

$ws.Range("D50").Value = "2.206.97"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("E51").Value = "  -2.71%  "
